# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Updates the IG metadata sheet (Version/Date bump, Publisher/Jurisdiction added in
# place of the duplicated "Contact" rows) and refreshes the root Extension's
# Short/Definition text on the Elements sheet.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Publication date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Row 9 stays "Publisher", but now carries a value
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a duplicated, unhelpful "Contact" row; replace with Jurisdiction
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second (duplicate) "Contact" / "No display for ContactDetail" row
# -- remove it entirely so everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# Update the root Extension element's Short/Definition on the Elements sheet
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Dual Eligibility Indicator"
$elements.Range("L2").Value = "Indicator of dual eligibility for the Medicaid market"
